# Update the "Förändrad" (Changed) date column (C) for rows 2-24
# from serial date 46066 (2026-02-13) to 46070 (2026-02-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46066) {
        $cell.Value2 = 46070
    }
}
